# Changed the tool statuses to numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$statusValues = @{
    2  = 4
    3  = 3
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 2
}

foreach ($row in $statusValues.Keys) {
    $ws.Range("J$row").Value = $statusValues[$row]
}

$ws.Range("J12").Select()
